$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 'The applicant has a diverse set of technical skills including JavaScript, PyTorch, Django, Data Analysis, Computer Vision, Python, and TensorFlow. Although the applicant''s projects do not directly mention experience with MongoDB or NodeJS, the projects showcase strong technical abilities in web development, backend development, AI/ML, and data analysis. The applicant''s experience in developing user-friendly platforms and working on cutting-edge AI/ML technologies demonstrates the potential to excel in a Remote MERN stack developer intern role.'
$ws.Range("E2").Value = 'The applicant has expressed a strong interest in working in Japan and has demonstrated adaptability and willingness to embrace a new culture. The applicant''s eagerness to learn a new language, appreciation of Japanese work culture, and openness to working in diverse locations showcase a positive attitude. Furthermore, the applicant''s collaborative problem-solving approach and the ability to bridge the gap between technical and creative aspects highlight strong teamwork and problem-solving skills. Based on the provided answers, the applicant exhibits strong potential for adapting to a new work environment in Japan, thus deserving a score of 4.'
$ws.Range("C3").Value = 'The applicant has worked on diverse projects that showcase proficiency in web development and real-time application development, aligning with the job description''s emphasis on web/mobile applications using MongoDB, ExpressJS, ReactJS, and NodeJS. The SHAMIYANA APP demonstrates skills in Flutter, Dart, and Firebase, while the RAPID project exhibits expertise in NodeJS, ExpressJS, Socket.IO, WebRTC, and Docker. Moreover, the SMART SENSING MIDDLEWARE project illustrates competence in ReactJS, HTML, and CSS, which are relevant to web development. Although the applicant lacks experience specifically with MongoDB, their strong foundation in web technologies and real-time applications positions them as a highly suitable candidate for the SDE Intern role.'
$ws.Range("E3").Value = 'The applicant has expressed a strong interest in Japan and its work culture, demonstrating adaptability and a willingness to embrace new experiences. Additionally, the applicant''s emphasis on team collaboration and problem-solving, as evidenced by their response to overcoming technical challenges, signifies a favorable personality fit for the company. Therefore, the applicant is given a score of 4 for demonstrating a positive attitude towards working in Japan and exhibiting key soft skills sought by the company.'
$ws.Range("C4").Value = 'The applicant has strong technical skills in Python, PyTorch, Tensorflow, and ReactJS which align well with the company''s requirements. Although the applicant lacks direct experience in MongoDB, ExpressJS, and NodeJS, the applicant''s expertise in web development and JavaScript, as well as experience with similar technologies, demonstrates the ability to quickly adapt and learn new technologies. The applicant''s projects showcase a blend of machine learning, web development, and AI, indicating a diverse skill set that can contribute to high-quality web/mobile applications. The applicant''s score reflects the alignment of their skills and project experiences with the job responsibilities and technical requirements.'
$ws.Range("E4").Value = 'The applicant has displayed strong communication skills by articulating a clear desire to be part of Japan''s work culture and expressing willingness to adapt to new cultures. Their commitment to learning new skills and the ability to work well in a team is evident from their responses. Furthermore, their demonstration of problem-solving abilities during challenging project situations signifies their potential to handle complexities in a new work environment. Overall, the applicant''s responses convey a positive attitude and a strong willingness to work in Japan, earning them a score of 4.'
